# Appends the new Gemini / corrected-Llama model-tournament pairings
# (rows 60-93) to Sheet1, and updates the active selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newPairs = @(
    @("gpt-4o", "gemini-1.5-flash-latest"),
    @("gpt-4o", "gemini-1.5-pro-latest"),
    @("gpt-4o-mini", "gemini-1.5-flash-latest"),
    @("gpt-4o-mini", "gemini-1.5-pro-latest"),
    @("gpt-3.5-turbo", "gemini-1.5-flash-latest"),
    @("gpt-3.5-turbo", "gemini-1.5-pro-latest"),
    @("grok-beta", "gemini-1.5-flash-latest"),
    @("grok-beta", "gemini-1.5-pro-latest"),
    @("claude-3-5-sonnet-latest", "gemini-1.5-flash-latest"),
    @("claude-3-5-sonnet-latest", "gemini-1.5-pro-latest"),
    @("claude-3-5-haiku-latest", "gemini-1.5-flash-latest"),
    @("claude-3-5-haiku-latest", "gemini-1.5-pro-latest"),
    @("gemini-1.5-flash-latest", "gpt-4o"),
    @("gemini-1.5-flash-latest", "gpt-4o-mini"),
    @("gemini-1.5-flash-latest", "gpt-3.5-turbo"),
    @("gemini-1.5-flash-latest", "grok-beta"),
    @("gemini-1.5-flash-latest", "claude-3-5-sonnet-latest"),
    @("gemini-1.5-flash-latest", "claude-3-5-haiku-latest"),
    @("gemini-1.5-flash-latest", "gemini-1.5-pro-latest"),
    @("gemini-1.5-flash-latest", "llama3.2-11b-vision"),
    @("gemini-1.5-flash-latest", "llama3.2-90b-vision"),
    @("gemini-1.5-pro-latest", "gpt-4o"),
    @("gemini-1.5-pro-latest", "gpt-4o-mini"),
    @("gemini-1.5-pro-latest", "gpt-3.5-turbo"),
    @("gemini-1.5-pro-latest", "grok-beta"),
    @("gemini-1.5-pro-latest", "claude-3-5-sonnet-latest"),
    @("gemini-1.5-pro-latest", "claude-3-5-haiku-latest"),
    @("gemini-1.5-pro-latest", "gemini-1.5-flash-latest"),
    @("gemini-1.5-pro-latest", "llama3.2-11b-vision"),
    @("gemini-1.5-pro-latest", "llama3.2-90b-vision"),
    @("llama-3.2-11b-vision", "gemini-1.5-flash-latest"),
    @("llama3.2-11b-vision", "gemini-1.5-pro-latest"),
    @("llama3.2-90b-vision", "gemini-1.5-flash-latest"),
    @("llama3.2-90b-vision", "gemini-1.5-pro-latest")
)

$startRow = 60
for ($i = 0; $i -lt $newPairs.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newPairs[$i][0]
    $ws.Cells.Item($r, 2).Value = $newPairs[$i][1]
}

# Matches the saved selection state recorded in the workbook.
[void]$ws.Range("G18").Select()
